# Power_ImpExpHubs.xlsx - add "Import Type" / "Export Type" columns (F:G)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power ImpExp")

# --- Column widths for F:G (new custom width used for the two new headers) ---
$ws.Range("F1:G1").ColumnWidth = 31.28515625

# --- Copy formatting from neighbouring columns onto the new F:G cells ---
$ws.Range("D3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("F4:G4").PasteSpecial(-4122)

$ws.Range("C5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("F6:G6").PasteSpecial(-4122)

$ws.Range("D7:E9").Copy()
$ws.Range("F7:G9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Updated Pmax Import / Pmax Export numbers ---
$ws.Range("D7").Value = 125
$ws.Range("E7").Value = 100
$ws.Range("D8").Value = 75
$ws.Range("E8").Value = 150
$ws.Range("D9").Value = 130
$ws.Range("E9").Value = 200

# --- New cell text, entered in the same order the author typed it ---
$ws.Range("F3").Value = "Import Type"
$ws.Range("G3").Value = "Export Type"

$ws.Range("F6").Value = "[ImpFix or ImpMax]"
$ws.Range("G6").Value = "[ExpFix or ExpMax]"

$ws.Range("F9").Value = "ImpFix"
$ws.Range("G9").Value = "ExpMax"

$ws.Range("F7").Value = "ImpMax"
$ws.Range("G7").Value = "ExpFix"
$ws.Range("F8").Value = "ImpMax"
$ws.Range("G8").Value = "ExpFix"

$ws.Range("F4").Value = "Determines wether ImpExp will be enforced as == ('Fix') or >= ('Max')"
$ws.Range("F5").Value = "Only ImpFix or ImpMax per hub, and ExpFix or ExpMax per hub"

# --- Dropdown data validation lists ---
$ws.Range("F7:F9").Validation.Add(3, 1, 1, '"ImpFix, ImpMax"')
$ws.Range("F7:F9").Validation.IgnoreBlank = $true
$ws.Range("F7:F9").Validation.InCellDropdown = $true
$ws.Range("F7:F9").Validation.ShowInput = $true
$ws.Range("F7:F9").Validation.ShowError = $true

$ws.Range("G7:G9").Validation.Add(3, 1, 1, '"ExpFix, ExpMax"')
$ws.Range("G7:G9").Validation.IgnoreBlank = $true
$ws.Range("G7:G9").Validation.InCellDropdown = $true
$ws.Range("G7:G9").Validation.ShowInput = $true
$ws.Range("G7:G9").Validation.ShowError = $true

Write-Host "done"
